$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("C7").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("B10").Value = "[-, Elcio Dec.-Desenho tecnico mecanico-1A]"
$ws.Range("D10").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("B11").Value = "[Ismail-Metrologia 1-1A, Andre B.-Comandos Eletricos-1A, -, -]"
$ws.Range("B12").Value = "[Ismail-Metrologia 1-1A, Andre B.-Comandos Eletricos-1A, -, -]"
$ws.Range("C12").Value = "[Aline S. M.-T. M. Metalicos-1A, Aline S. M.-T. M. Metalicos-1A]"
$ws.Range("B14").Value = "[Ismail-Metrologia 1-1A, Andre B.-Comandos Eletricos-1A, -, -]"
$ws.Range("C14").Value = "Anselmo-Gestao Integrada"
$ws.Range("B15").Value = "[Ismail-Metrologia 1-1A, Andre B.-Comandos Eletricos-1A, -, -]"
$ws.Range("C15").Value = "Anselmo-Gestao Integrada"
$ws.Range("B16").Value = "[Elcio Dec.-Desenho tecnico mecanico-1A, Elcio Dec.-Desenho tecnico mecanico-1A]"
$ws.Range("D18").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("C20").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("D20").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("C21").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
$ws.Range("E21").Value = "[Ismail-Metrologia 1-1A, -, -, -]"
